$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = 1.63

$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("U3").Value = 2.37
$ws.Range("V3").Value = 1.5

$ws.Range("V4").Value = 1.54

$ws.Range("V5").Value = 1.58

$ws.Range("V6").Value = 1.54

$ws.Range("V7").Value = 1.58

$ws.Range("V8").Value = 1.72
